$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on D/E columns touched below so purely-numeric-looking
# strings (e.g. "577.70", "1.00") are not auto-coerced to Number by the Value setter,
# matching the inlineStr/text storage used throughout this sheet.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.540.62"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.30%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.459.01"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.78%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.70"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.52%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.18"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.47%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.457.44"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.61%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.556"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.18%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.57"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.72%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.55%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.89%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.058.12"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.87%  "

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.90%  "

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +9.89%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.85"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.94%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.565.46"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.17%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.466.12"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.86%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.44"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.17%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.36"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.91%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "396.17"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.62%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.52"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.93%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.547"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.61%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.03"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.39%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.27%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000122"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +27.28%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +7.35%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.60%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.27%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.08"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +9.71%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.50%  "

# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.39"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +7.32%  "

# Row 33
$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.68"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.34%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.84"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.96%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.04%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.07"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.23%  "

# Row 37
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.59"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.51%  "

# Row 38
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.47"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.77%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0785"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.10%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.88"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.59%  "

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.87%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.938.26"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.75%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.60%  "

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.35%  "

# Row 45
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.44"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.52%  "

# Row 46
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.24"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.27%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.92"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +9.12%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +25.52%  "

# Row 50
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.861"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +7.32%  "

# Row 51
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.56"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.33%  "
